# Updated cryptos list on Wed Apr 26 08:52:44 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.705.35"
$ws.Range("E2").Value = "  +4.51%  "
$ws.Range("D3").Value = "1.874.20"
$ws.Range("E3").Value = "  +2.64%  "
$ws.Range("D4").Value = "'0.9985"
$ws.Range("E4").Value = "  -0.74%  "
$ws.Range("D5").Value = "'338.86"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").Value = "'0.9992"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").Value = "'0.4707"
$ws.Range("E7").Value = "  +3.30%  "
$ws.Range("D8").Value = "'0.4002"
$ws.Range("E8").Value = "  +5.17%  "
$ws.Range("D9").Value = "'47.71"
$ws.Range("E9").Value = "  +2.84%  "
$ws.Range("D10").Value = "'0.08058"
$ws.Range("E10").Value = "  +2.29%  "
$ws.Range("D11").Value = "'1.006"
$ws.Range("E11").Value = "  +3.52%  "
$ws.Range("D12").Value = "'22.17"
$ws.Range("E12").Value = "  +5.63%  "
$ws.Range("D13").Value = "'6.066"
$ws.Range("E13").Value = "  +3.43%  "
$ws.Range("D14").Value = "1.863.13"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").Value = "'7.294"
$ws.Range("E15").Value = "  +3.77%  "
$ws.Range("D16").Value = "'90.49"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").Value = "'0.9992"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "'0.00001045"
$ws.Range("E18").Value = "  +1.85%  "
$ws.Range("D19").Value = "'0.06612"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "'17.67"
$ws.Range("E20").Value = "  +2.94%  "
$ws.Range("D21").Value = "'0.9993"
$ws.Range("D22").Value = "28.688.91"
$ws.Range("E22").Value = "  +4.51%  "
$ws.Range("D23").Value = "'5.513"
$ws.Range("E23").Value = "  +3.54%  "
$ws.Range("D24").Value = "'11.07"
$ws.Range("E24").Value = "  +2.56%  "
$ws.Range("D25").Value = "'2.259"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").Value = "2.084.98"
$ws.Range("E26").Value = "  +2.06%  "
$ws.Range("D27").Value = "'160.52"
$ws.Range("E27").Value = "  +2.38%  "
$ws.Range("D28").Value = "'19.81"
$ws.Range("E28").Value = "  +2.08%  "
$ws.Range("D29").Value = "'2.134"
$ws.Range("E29").Value = "  +3.59%  "
$ws.Range("D30").Value = "'5.508"
$ws.Range("E30").Value = "  +4.99%  "
$ws.Range("D31").Value = "'120.01"
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("D32").Value = "'0.9818"
$ws.Range("E32").Value = "  +3.72%  "
$ws.Range("D33").Value = "'0.09541"
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("D34").Value = "'3.675"
$ws.Range("E34").Value = "  +2.64%  "
$ws.Range("D35").Value = "'1.392"
$ws.Range("E35").Value = "  +5.39%  "
$ws.Range("D36").Value = "'5.381"
$ws.Range("E36").Value = "  +2.80%  "
$ws.Range("D37").Value = "'0.06210"
$ws.Range("E37").Value = "  +4.86%  "
$ws.Range("D38").Value = "'0.02263"
$ws.Range("E38").Value = "  +3.83%  "
$ws.Range("D39").Value = "'8.475"
$ws.Range("E39").Value = "  +5.66%  "
$ws.Range("D40").Value = "'1.182"
$ws.Range("E40").Value = "  +1.83%  "
$ws.Range("D41").Value = "'0.5968"
$ws.Range("E41").Value = "  +3.93%  "
$ws.Range("D42").Value = "'0.9990"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("E43").Value = "  +3.13%  "
$ws.Range("D44").Value = "'10.36"
$ws.Range("E44").Value = "  +3.51%  "
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("D46").Value = "'0.5590"
$ws.Range("E46").Value = "  +2.67%  "

# Row 47/48: EnergySwap overtakes Cronos in ranking, so the two rows swap places
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'12.19"
$ws.Range("E47").Value = "  +1.17%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.07389"
$ws.Range("E48").Value = "  +11.77%  "

$ws.Range("D49").Value = "'1.967"
$ws.Range("E49").Value = "  +5.56%  "
$ws.Range("D50").Value = "'2.083"
$ws.Range("E50").Value = "  +12.69%  "
$ws.Range("D51").Value = "'112.61"
$ws.Range("E51").Value = "  +2.10%  "
